# New PO forecast model
# Updates the three PO data sheets ("Weekly Quantity", "Monthly Trend",
# "PO Forecast") with the refreshed forecast numbers and appends the
# latest data point to each sheet.

$wb = $excel.ActiveWorkbook

$dateFormat = "YYYY-MM-DD HH:MM:SS"

# ---------------------------------------------------------------------
# Sheet: Weekly Quantity  -> append a new weekly data point (row 10)
# ---------------------------------------------------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Cells.Item(10, 1).Value = 45669.99999999999
$wsWeekly.Cells.Item(10, 1).NumberFormat = $dateFormat
$wsWeekly.Cells.Item(10, 2).Value = 1

# ---------------------------------------------------------------------
# Sheet: Monthly Trend -> append a new monthly data point (row 9)
# ---------------------------------------------------------------------
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Cells.Item(9, 1).Value = 45688.99999999999
$wsMonthly.Cells.Item(9, 1).NumberFormat = $dateFormat
$wsMonthly.Cells.Item(9, 2).Value = 1

# ---------------------------------------------------------------------
# Sheet: PO Forecast -> refreshed forecast model values
# ---------------------------------------------------------------------
$wsForecast = $wb.Worksheets.Item("PO Forecast")

# Quantities for the already-existing historical weeks (rows 2-9) change
# under the new model, while the dates stay the same.
$wsForecast.Cells.Item(2, 2).Value = 27
$wsForecast.Cells.Item(3, 2).Value = 1
$wsForecast.Cells.Item(4, 2).Value = 1
$wsForecast.Cells.Item(5, 2).Value = 1
$wsForecast.Cells.Item(6, 2).Value = 3
$wsForecast.Cells.Item(7, 2).Value = 2
$wsForecast.Cells.Item(8, 2).Value = 2
$wsForecast.Cells.Item(9, 2).Value = 2

# Rows 10-17 are replaced with the new forecast horizon (new dates and
# quantities), and a brand-new row 18 is appended.
$wsForecast.Cells.Item(10, 1).Value = 45669.99999999999
$wsForecast.Cells.Item(10, 1).NumberFormat = $dateFormat
$wsForecast.Cells.Item(10, 2).Value = 1

$wsForecast.Cells.Item(11, 1).Value = 45676.99999999999
$wsForecast.Cells.Item(11, 1).NumberFormat = $dateFormat
$wsForecast.Cells.Item(11, 2).Value = 139

$wsForecast.Cells.Item(12, 1).Value = 45683.99999999999
$wsForecast.Cells.Item(12, 1).NumberFormat = $dateFormat
$wsForecast.Cells.Item(12, 2).Value = 0

$wsForecast.Cells.Item(13, 1).Value = 45690.99999999999
$wsForecast.Cells.Item(13, 1).NumberFormat = $dateFormat
$wsForecast.Cells.Item(13, 2).Value = 0

$wsForecast.Cells.Item(14, 1).Value = 45697.99999999999
$wsForecast.Cells.Item(14, 1).NumberFormat = $dateFormat
$wsForecast.Cells.Item(14, 2).Value = 0

$wsForecast.Cells.Item(15, 1).Value = 45704.99999999999
$wsForecast.Cells.Item(15, 1).NumberFormat = $dateFormat
$wsForecast.Cells.Item(15, 2).Value = 0

$wsForecast.Cells.Item(16, 1).Value = 45711.99999999999
$wsForecast.Cells.Item(16, 1).NumberFormat = $dateFormat
$wsForecast.Cells.Item(16, 2).Value = 41

$wsForecast.Cells.Item(17, 1).Value = 45718.99999999999
$wsForecast.Cells.Item(17, 1).NumberFormat = $dateFormat
$wsForecast.Cells.Item(17, 2).Value = 36

$wsForecast.Cells.Item(18, 1).Value = 45725.99999999999
$wsForecast.Cells.Item(18, 1).NumberFormat = $dateFormat
$wsForecast.Cells.Item(18, 2).Value = 0
